$d = $word.ActiveDocument

# Locate the end of the first paragraph's text (after "group.") and split
# it into a new paragraph, then type the new sentence.
$rng = $d.Content
$rng.Find.Execute("exploratory group.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertParagraphAfter()

# Move into the newly created paragraph (the second paragraph) and add text
$rng2 = $d.Paragraphs(2).Range
$rng2.InsertBefore("We can offer some excellent activates.")

# Now add a run with two spaces after the bookmark at the end of paragraph 2
$endRng = $d.Paragraphs(2).Range
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1) | Out-Null
